# NIT-9000049211.xlsx — "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta".
#
# The worker table on the sheet used to list two workers:
#   row 16 -> CC 1143326345 / INGRID PAOLA BATISTA MUÑOZ / periodo 2501
#   row 17 -> CC 1010068328 / ISAIAS JAVIER BALMACEDA PRINS / periodo 2002
# The update drops Ingrid's row from the statement and keeps only Isaias'
# record (now the sole worker row), adjusting the header totals
# (Valor Mora, Cant. Trabajadores, Cant. Periodos) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Isaias' original row; the remaining worker row (formerly Ingrid's)
# keeps its place in the table and gets overwritten with Isaias' data below.
$ws.Rows.Item(17).Delete()

# Worker row (now row 16): CC / id / name / periodo / valor mora / salario basico
$ws.Cells.Item(16, 3).Value = "1010068328"
$ws.Cells.Item(16, 4).Value = "ISAIAS JAVIER BALMACEDA PRINS"
$ws.Cells.Item(16, 5).Value = "2002"
$ws.Cells.Item(16, 6).Value = 56999
$ws.Cells.Item(16, 7).Value = 1424973

# Header totals: now a single worker over a single period.
$ws.Cells.Item(11, 5).Value = 56999
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 6).Value = 1
